# "600 Essential Words.xlsx" - add vocabulary rows for determine / engage /
# establish / obligate / provision / resolve / specific on the "Contracts"
# sheet, fill in the missing IPA for "abide by", widen column G, and move
# the active sheet/selection from "Bussines Planning" back to "Contracts".

$wb = $excel.ActiveWorkbook
$wsContracts = $wb.Worksheets.Item("Contracts")
$wsPlanning  = $wb.Worksheets.Item("Bussines Planning")

function Set-PlainCell($ws, $cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

function Set-RichCell($ws, $cellRef, $boldText, $restText) {
    # Whole string first (plain), then re-apply bold + matching font to the
    # leading run so the shared-string gets split into two <r> runs, same
    # as Excel does when you bold part of a cell's text.
    $full = $boldText + $restText
    $ws.Range($cellRef).Value = $full

    $boldLen = $boldText.Length
    $restLen = $restText.Length

    $ws.Range($cellRef).Characters(1, $boldLen).Font.Bold = $true
    $ws.Range($cellRef).Characters(1, $boldLen).Font.Size = 13
    $ws.Range($cellRef).Characters(1, $boldLen).Font.Name = "Times New Roman"

    $ws.Range($cellRef).Characters($boldLen + 1, $restLen).Font.Size = 13
    $ws.Range($cellRef).Characters($boldLen + 1, $restLen).Font.Name = "Times New Roman"
}

# --- row 6: determine ------------------------------------------------------
Set-PlainCell $wsContracts "B6" "v"
Set-PlainCell $wsContracts "G6" "After reading the contract, I was still unable to determine if our company was liable for back wages."
Set-PlainCell $wsContracts "D6" "/dɪˈtɜːmɪn/"
Set-PlainCell $wsContracts "F6" "to discover the facts about something; to calculate something exactly"
Set-PlainCell $wsContracts "C6" "to find out; to influence; establish"
Set-PlainCell $wsContracts "E6" "Quyết định, xác định"

# --- row 7: engage -----------------------------------------------------------
Set-PlainCell $wsContracts "B7" "v"
Set-PlainCell $wsContracts "D7" "/ɪnˈɡeɪdʒ/"
Set-RichCell  $wsContracts "F7" "engage something/somebody (formal)" " to succeed in attracting and keeping somebody’s attention and interest"
Set-PlainCell $wsContracts "C7" "participate; involve"
Set-PlainCell $wsContracts "E7" "Thu hút"
Set-PlainCell $wsContracts "G7" "Before engage in a new business, it is important to do thorough research."

# --- row 8: establish --------------------------------------------------------
Set-PlainCell $wsContracts "B8" "v"
Set-PlainCell $wsContracts "G8" "Through her many books and interview, Dr.Wan established herself as an authority on conflict resolution."
Set-PlainCell $wsContracts "D8" "/ɪˈstæblɪʃ/"
Set-PlainCell $wsContracts "C8" "to institule permanently; to bring about, set up"
Set-RichCell  $wsContracts "F8" "establish something" " to start or create an organization, a system, etc. that is meant to last for a long time"
Set-PlainCell $wsContracts "E8" "Thành lập, chứng minh, củng cố"

# --- fill the previously-blank IPA cell for "abide by" ------------------------
Set-PlainCell $wsContracts "D2" "/əˈbaɪd/"

# --- row 9: obligate ----------------------------------------------------------
Set-PlainCell $wsContracts "B9" "v"
Set-PlainCell $wsContracts "C9" "to blind legally ỏ morally"
Set-PlainCell $wsContracts "G9" "The contractor was obligated by the contract to work 40 hours a week."
Set-PlainCell $wsContracts "E9" "Bắt buộc."
Set-PlainCell $wsContracts "D9" "/ˈɒb.lɪ.ɡeɪt/"
Set-PlainCell $wsContracts "F9" "to force someone to do something, or to make it necessary for someone to do something"

# --- row 11: provision ---------------------------------------------------------
Set-PlainCell $wsContracts "B11" "n"
Set-PlainCell $wsContracts "D11" "/prəˈvɪʒn/"
Set-PlainCell $wsContracts "G11" "The father made provisions for his children through his will."
Set-RichCell  $wsContracts "F11" "provision for somebody/something" " preparations that you make for something that might or will happen in the future"
Set-PlainCell $wsContracts "E11" "cung cấp , giao kèo điều khoản"
Set-PlainCell $wsContracts "C11" "a measure  taken before;a stipulation"

# --- row 12: resolve -----------------------------------------------------------
Set-PlainCell $wsContracts "B12" "v"
Set-PlainCell $wsContracts "D12" "/rɪˈzɒlv/"
Set-PlainCell $wsContracts "F12" "to find an acceptable solution to a problem or difficulty"
Set-PlainCell $wsContracts "G12" "The mediator was able to resolve the problem to everyone's satisfaction."
Set-PlainCell $wsContracts "C12" "to deal with successfully; to declare; conviction"
Set-PlainCell $wsContracts "E12" "kiên quyết, thông qua"

# --- row 13: specific -----------------------------------------------------------
Set-PlainCell $wsContracts "B13" "adj"
Set-PlainCell $wsContracts "E13" "rõ ràng, cụ thể, đặc trung riêng biệt"
Set-PlainCell $wsContracts "D13" "/spəˈsɪfɪk/"
Set-PlainCell $wsContracts "F13" " connected with one particular thing only"
Set-PlainCell $wsContracts "C13" "particular, precise, peculiar"
Set-PlainCell $wsContracts "G13" "The customer's specific complanit not addressed in his e-mail."

# --- widen column G (Original Sentencs) to fit the new long sentences ----------
$wsContracts.Columns("G:G").ColumnWidth = 98.85546875

# --- move the selection/active tab: "Bussines Planning" loses focus, -----------
# "Contracts" becomes the active sheet/tab again.
$wsPlanning.Activate()
$wsPlanning.Range("B37").Select()

$wsContracts.Activate()
$wsContracts.Range("G19").Select()
